# "import excel to update kanban"
# The kanban template header row lost its "Wire Nr" column (old column F) and
# the trailing "Operate" header was corrected to "Operator". Also refresh the
# page setup (paper size / orientation) that came along with the re-import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Wire Nr" column entirely - this shifts every column
# after it (Product Nr, Type, Bundle, ... Operate) one slot to the left,
# which is exactly what the sharedStrings/sheet diff shows.
$ws.Columns("F").Delete()

# Fix the trailing header's wording.
$ws.Range("N1").Value = "Operator"

# Make the freshly edited cell the active selection, as Excel would leave it
# right after typing the correction.
[void]$ws.Range("N1").Select()

# Page setup that accompanied the re-import of this sheet.
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
